# Updates all PDF figures: refresh the underlying case-count inputs
# (columns B:E, rows 19-38) on the measles_costs sheet, then realign
# the window/selection view state to match the author's last save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("measles_costs")

# New case-count figures (Cases, Doses, Cases PCR confirmed, Contacts) per outbreak row.
$data = @{
    19 = @(436350, 74410, 10662, 40320)
    20 = @(206000, 35150, 5041, 19056)
    21 = @(482180, 82195, 11770, 44525)
    22 = @(283700, 48368, 6928, 26204)
    23 = @(469300, 80484, 11633, 43820)
    24 = @(151700, 25934, 3730, 14082)
    25 = @(138380, 23864, 3478, 13053)
    26 = @(98196, 16854, 2439, 9182)
    27 = @(162560, 27397, 3855, 14697)
    28 = @(137000, 23502, 3398, 12799)
    29 = @(151690, 25902, 3719, 14051)
    30 = @(55620, 9444, 1344, 5099)
    31 = @(297420, 50047, 7025, 26811)
    32 = @(43650, 7473, 1077, 4063)
    33 = @(109750, 18842, 2728, 10268)
    34 = @(359310, 61138, 8730, 33067)
    35 = @(41112, 7005, 1002, 3793)
    36 = @(525550, 90522, 13170, 49463)
    37 = @(32151, 5554, 812, 3042)
    38 = @(60120, 10218, 1457, 5521)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    # The pasted-in figures carry no number formatting (unlike the old s="2" cells).
    $rng = $ws.Range($ws.Cells.Item($row, 2), $ws.Cells.Item($row, 5))
    $rng.Style = "Normal"
}

# Recalculate so every dependent formula cell carries a fresh cached value.
$excel.Calculate()

# Restore the view state captured at the author's last save.
$ws.Range("B14").Select()

$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$win.Left = 14505
$win.Top = -15
$win.Width = 14310
$win.Height = 12855
